$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy header formatting (bold, border, centered) from H1 into the new I1/J1 header cells
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

$data = @(
  @(1, 2),
  @(1, 6),
  @(1, 6),
  @(1, 4),
  @(1, 6),
  @(1, 3),
  @(1, 6),
  @(1, 5),
  @(6, 8),
  @(7, 7),
  @(6, 8),
  @(7, 7),
  @(1, 2),
  @(5, 7),
  @(4, 4),
  @(10, 10),
  @(9, 9),
  @(8, 8),
  @(5, 6),
  @(5, 7),
  @(1, 3),
  @(1, 2)
)

for ($i = 0; $i -lt $data.Length; $i++) {
  $row = $i + 2
  $ws.Cells.Item($row, 9).Value = $data[$i][0]
  $ws.Cells.Item($row, 10).Value = $data[$i][1]
}
